# Fix the #DIV/0! errors ("null pointer exception" in the commit message)
# on RQ1 and RQ4 by filling in the missing numerator/denominator values
# that were left blank, and update the view/selection state to match.

$wb = $excel.ActiveWorkbook

# --- RQ1 (sheet1): row 6 ("Carol") was missing C6/D6/F6/G6, causing
#     E6 (=C6/D6) and H6 (=F6/G6) to evaluate to #DIV/0!. Row 7's F7 also
#     changes from 8 to 10.
$ws1 = $wb.Worksheets.Item("RQ1")
$ws1.Range("C6").Value = 298
$ws1.Range("D6").Value = 24
$ws1.Range("F6").Value = 22
$ws1.Range("G6").Value = 24
$ws1.Range("F7").Value = 10

# --- RQ4 (sheet4): row 7 ("Jabref") was missing B7/C7/E7/F7, causing
#     D7 (=C7/B7*100) and G7 (=F7/E7*100) to evaluate to #DIV/0!.
$ws4 = $wb.Worksheets.Item("RQ4")
$ws4.Range("B7").Value = 35
$ws4.Range("C7").Value = 7
$ws4.Range("E7").Value = 3
$ws4.Range("F7").Value = 2

# --- Update selections / active sheet to match the saved view state.
$ws1.Range("F6").Select() | Out-Null

$ws3 = $wb.Worksheets.Item("RQ3")
$ws3.Range("E7").Select() | Out-Null

$ws4.Activate() | Out-Null
$ws4.Range("F7").Select() | Out-Null
